$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (D) and volume-change (E) values to match latest scrape.
# D-column values are plain numeric-looking strings (e.g. "321.55", "48.253.49")
# that must stay as literal text (matching the source inlineStr cells), so we
# force a temporary Text number format before assigning, then restore the default
# "Normal" style so no stray formatting is introduced.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '48.253.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.508.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.527'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.01%  '
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.901.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.506.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '48.101.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.97%  '
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0947'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '276.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E28').Value = '  +4.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.46%  '
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0784'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.22'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0306'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.000.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.12%  '
$ws.Range('E47').Value = '  +3.77%  '
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.51%  '
